# ---------------------------------------------------------------------------
# Adds the "Completado Por" column + sample rows to the tasks sheet, renames
# the sheet to "Tareas", and adds a new "Instrucciones" help sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet and add the new instructions sheet right
#    after it.
# ---------------------------------------------------------------------------
$tareas = $wb.Worksheets.Item(1)
$tareas.Name = "Tareas"

$instrucciones = $wb.Worksheets.Add($null, $tareas)
$instrucciones.Name = "Instrucciones"

# ---------------------------------------------------------------------------
# 2. "Tareas" sheet: widen/add columns, write header + sample rows.
# ---------------------------------------------------------------------------

# Column widths (values chosen so the saved width matches the target widths
# exactly, accounting for the character->pixel rounding Excel applies).
$tareas.Columns.Item(2).ColumnWidth = 44.083333333333336   # B -> 45
$tareas.Columns.Item(5).ColumnWidth = 24.083333333333332   # E -> 25
$tareas.Columns.Item(6).ColumnWidth = 29.083333333333332   # F -> 30
$tareas.Columns.Item(7).ColumnWidth = 13.083333333333334   # G -> 14
$tareas.Columns.Item(8).ColumnWidth = 11.083333333333334   # H -> 12
$tareas.Columns.Item(9).ColumnWidth = 17.083333333333332   # I -> 18 (new)

# Force all the data cells to Text format first so values that look like
# dates/numbers (e.g. "2024-12-20") are NOT auto-converted by the engine.
# G3/G4 are overridden back to real numbers further below.
$tareas.Range("A2:F4").NumberFormat = "@"
$tareas.Range("H2:I4").NumberFormat = "@"

# Header row
$tareas.Cells.Item(1,1).Value = "Titulo"
$tareas.Cells.Item(1,2).Value = "Descripcion"
$tareas.Cells.Item(1,3).Value = "Prioridad"
$tareas.Cells.Item(1,4).Value = "Fecha Vencimiento"
$tareas.Cells.Item(1,5).Value = "Asignados"
$tareas.Cells.Item(1,6).Value = "Etiquetas"
$tareas.Cells.Item(1,7).Value = "Tiempo (min)"
$tareas.Cells.Item(1,8).Value = "Estado"
$tareas.Cells.Item(1,9).Value = "Completado Por"

# Header style: bold white text, blue fill, thin border, centered.
$header = $tareas.Range("A1:I1")
$header.Font.Bold = $true
$header.Font.Color = 16777215
$header.Interior.Color = 15426341
$header.Interior.PatternColor = 15426341
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108

# Row 2 (replaces old sample task)
$tareas.Cells.Item(2,1).Value = "Revisar expediente 1234"
$tareas.Cells.Item(2,2).Value = "Verificar documentacion completa del caso"
$tareas.Cells.Item(2,3).Value = "Normal"
$tareas.Cells.Item(2,4).Value = "2024-12-20"
$tareas.Cells.Item(2,5).Value = "admin"
$tareas.Cells.Item(2,6).Value = "Urgente, Legal"
$tareas.Cells.Item(2,7).NumberFormat = "@"
$tareas.Cells.Item(2,7).Value = ""
$tareas.Cells.Item(2,8).Value = "Pendiente"
$tareas.Cells.Item(2,9).NumberFormat = "@"
$tareas.Cells.Item(2,9).Value = ""

# Row 3 (new)
$tareas.Cells.Item(3,1).Value = "Preparar informe mensual"
$tareas.Cells.Item(3,2).Value = "Elaborar informe de actividades del mes"
$tareas.Cells.Item(3,3).Value = "Media"
$tareas.Cells.Item(3,4).Value = "2024-12-15"
$tareas.Cells.Item(3,5).Value = "admin"
$tareas.Cells.Item(3,6).Value = "Administrativo"
$tareas.Cells.Item(3,7).Value = 120
$tareas.Cells.Item(3,8).Value = "Completada"
$tareas.Cells.Item(3,9).Value = "giuliana"

# Row 4 (new)
$tareas.Cells.Item(4,1).Value = "Audiencia caso Smith"
$tareas.Cells.Item(4,2).Value = "Preparar alegatos para audiencia"
$tareas.Cells.Item(4,3).Value = "Urgente"
$tareas.Cells.Item(4,4).Value = "2024-12-10"
$tareas.Cells.Item(4,5).Value = "admin, giuliana"
$tareas.Cells.Item(4,6).Value = "Legal, Tribunal"
$tareas.Cells.Item(4,7).Value = 60
$tareas.Cells.Item(4,8).Value = "Completada"
$tareas.Cells.Item(4,9).NumberFormat = "@"
$tareas.Cells.Item(4,9).Value = ""

# Thin border around all the data cells (rows 2-4, columns A-I)
$tareas.Range("A2:I4").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3. "Instrucciones" sheet: title + numbered help text.
# ---------------------------------------------------------------------------
$instrucciones.Columns.Item(1).ColumnWidth = 74.08333333333333   # A -> 75

$instrucciones.Cells.Item(1,1).Value = "INSTRUCCIONES PARA IMPORTAR TAREAS"
$instrucciones.Cells.Item(1,1).Font.Bold = $true
$instrucciones.Cells.Item(1,1).Font.Size = 14

$instrucciones.Cells.Item(3,1).Value = "1. Complete la hoja ""Tareas"" con los datos de las tareas a importar."

$instrucciones.Cells.Item(5,1).Value = "2. Columnas obligatorias:"
$instrucciones.Cells.Item(6,1).Value = "   - Titulo: Nombre de la tarea"
$instrucciones.Cells.Item(7,1).Value = "   - Prioridad: Normal, Media o Urgente"
$instrucciones.Cells.Item(8,1).Value = "   - Fecha Vencimiento: Formato AAAA-MM-DD, ej: 2024-12-20"
$instrucciones.Cells.Item(9,1).Value = "   - Asignados: Username(s) separados por coma, ej: admin, usuario1"

$instrucciones.Cells.Item(11,1).Value = "3. Columnas opcionales:"
$instrucciones.Cells.Item(12,1).Value = "   - Descripcion: Detalle de la tarea"
$instrucciones.Cells.Item(13,1).Value = "   - Etiquetas: Nombre(s) de etiquetas separadas por coma"
$instrucciones.Cells.Item(14,1).Value = "   - Tiempo (min): Tiempo dedicado en minutos, ej: 60, 120"
$instrucciones.Cells.Item(15,1).Value = "   - Estado: Pendiente o Completada (por defecto: Pendiente)"
$instrucciones.Cells.Item(16,1).Value = "   - Completado Por: Username del usuario que completo la tarea"
$instrucciones.Cells.Item(17,1).Value = "     * Si el estado es Completada y no se indica usuario, se usa quien sube el archivo"
$instrucciones.Cells.Item(18,1).Value = "     * Si el estado es Pendiente, este campo se ignora"

$instrucciones.Cells.Item(20,1).Value = "4. Importante:"
$instrucciones.Cells.Item(21,1).Value = "   - No modifique los encabezados de las columnas"
$instrucciones.Cells.Item(22,1).Value = "   - Los usernames deben existir en el sistema"
$instrucciones.Cells.Item(23,1).Value = "   - Las etiquetas deben existir en el sistema (se ignoran las que no existan)"
$instrucciones.Cells.Item(24,1).Value = "   - La fecha debe estar en formato AAAA-MM-DD"

# Italic style for the indented bullet rows.
$bulletRows = @(6,7,8,9,12,13,14,15,16,21,22,23,24)
foreach ($r in $bulletRows) {
    $instrucciones.Cells.Item($r,1).Font.Italic = $true
}

# ---------------------------------------------------------------------------
# 4. Leave the "Tareas" sheet selected/active (matches the original file).
# ---------------------------------------------------------------------------
$tareas.Activate()
